$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 6.140625

# --- G3:G9 date values with date-formatted style matching existing body fill ---
$ws.Range("G3:G9").Value = 44340
$ws.Range("B3").Copy()
$ws.Range("G3:G9").PasteSpecial(-4122)
$ws.Range("G3:G9").NumberFormat = "d-mmm"

# --- Row 10: Total row ---
# Plain cells (B10, C10, D10, F10, G10): fill matches header fill, no alignment, no bold
$plainCells = @("B10","C10","D10","F10","G10")
foreach ($addr in $plainCells) {
  $ws.Range("B2").Copy()
  $ws.Range($addr).PasteSpecial(-4122)
  $ws.Range($addr).Font.Bold = $false
  $ws.Range($addr).HorizontalAlignment = 1
  $ws.Range($addr).VerticalAlignment = -4107
}
$ws.Range("B10").Value = "Total"

# E10: total sum, with header fill, centered alignment, no bold
$ws.Range("B2").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Font.Bold = $false
$ws.Range("E10").Formula = "=SUM(E3:E9)"

# --- Sheet view / selection ---
$ws.Range("E10").Select()
